{"js": "// Word's auto-generated \"_GoBack\" bookmark originally sits at the end of the\n// title paragraph (\"... Sucursal\"); the image in the third paragraph is\n// removed and the \"_GoBack\" bookmark ends up wrapping the (now empty) third\n// paragraph instead.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\n\nconst inlinePictures = body.inlinePictures;\ninlinePictures.load(\"items\");\n\nawait context.sync();\n\n// Remove the picture that lives in the third paragraph.\nif (inlinePictures.items.length > 0) {\n  inlinePictures.items[0].delete();\n}\n\n// Drop the existing \"_GoBack\" bookmark (currently after \"Sucursal\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-create \"_GoBack\" around the third paragraph, which is now empty.\nconst targetParagraph = paragraphs.items[2];\nconst targetRange = targetParagraph.getRange();\ntargetRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word's auto-generated \"_GoBack\" bookmark originally sits at the end of the\n# title paragraph (\"... Sucursal\"); the image in the third paragraph is\n# removed and the \"_GoBack\" bookmark ends up wrapping the (now empty) third\n# paragraph instead.\n\n$d = $word.ActiveDocument\n\n# Remove the picture that lives in the third paragraph.\nif ($d.InlineShapes.Count -gt 0) {\n    $d.InlineShapes.Item(1).Delete()\n}\n\n# Drop the existing \"_GoBack\" bookmark (currently after \"Sucursal\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-create \"_GoBack\" around the third paragraph, which is now empty.\n$targetParagraph = $d.Paragraphs.Item(3)\n$d.Bookmarks.Add(\"_GoBack\", $targetParagraph.Range)\n"}
